$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "date"

$ws.Range("E2").Value = (Get-Date -Year 2020 -Month 2 -Day 21 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E3").Value = (Get-Date -Year 2020 -Month 2 -Day 21 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E4").Value = (Get-Date -Year 2020 -Month 2 -Day 25 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E5").Value = (Get-Date -Year 2020 -Month 4 -Day 20 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("E2:E4").NumberFormat = "yyyy-mm-dd;@"
$ws.Range("E2:E4").Borders.LineStyle = 1
$ws.Range("E2:E4").Font.Size = 9
$ws.Range("E2:E4").HorizontalAlignment = -4131
$ws.Range("E2:E4").VerticalAlignment = -4160
$ws.Range("E2:E4").ShrinkToFit = $true

$ws.Range("E5").NumberFormat = "m/d/yyyy"

$ws.Range("F7").Select()
